$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 791.5
$ws.Range("I18").Value = 790.7273
$ws.Range("J18").Value = 800
$ws.Range("K18").Value = 790.7273
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = -506.7273
$ws.Range("N18").Value = -1368

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2689.6924
$ws.Range("I100").Value = 1790
$ws.Range("J100").Value = 4129.2
$ws.Range("K100").Value = 1790
$ws.Range("L100").Value = 4129.2
$ws.Range("M100").Value = -1249
$ws.Range("N100").Value = -5211.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 42880.5
$ws.Range("J120").Value = 42880.5
$ws.Range("L120").Value = 42880.5
$ws.Range("N120").Value = -52556.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2208.4875
$ws.Range("I132").Value = 1801.8448
$ws.Range("J132").Value = 3280.5454
$ws.Range("K132").Value = 5405.5344
$ws.Range("L132").Value = 9841.636200000001
$ws.Range("M132").Value = -2875.5344
$ws.Range("N132").Value = -14901.6362

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 343.31033
$ws.Range("I135").Value = 341.8909
$ws.Range("J135").Value = 369.33334
$ws.Range("K135").Value = 3077.0181
$ws.Range("L135").Value = 3324.00006
$ws.Range("M135").Value = -542.0180999999998
$ws.Range("N135").Value = -8394.00006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2291.3804
$ws.Range("I137").Value = 877.375
$ws.Range("J137").Value = 5523.393
$ws.Range("K137").Value = 2632.125
$ws.Range("L137").Value = 16570.179
$ws.Range("M137").Value = -82.125
$ws.Range("N137").Value = -21670.179

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1742.83
$ws.Range("I32").Value = 1475.5591
$ws.Range("J32").Value = 5293.7144
$ws.Range("K32").Value = 1475.5591
$ws.Range("L32").Value = 5293.7144
$ws.Range("M32").Value = -1188.5591
$ws.Range("N32").Value = -5867.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 885.83826
$ws.Range("I61").Value = 689.0172
$ws.Range("J61").Value = 2027.4
$ws.Range("K61").Value = 689.0172
$ws.Range("L61").Value = 2027.4
$ws.Range("M61").Value = -477.0172
$ws.Range("N61").Value = -2451.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2290.2
$ws.Range("I74").Value = 450.53192
$ws.Range("J74").Value = 13098.25
$ws.Range("K74").Value = 450.53192
$ws.Range("L74").Value = 13098.25
$ws.Range("M74").Value = 423.46808
$ws.Range("N74").Value = -14846.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2290.2
$ws.Range("I77").Value = 450.53192
$ws.Range("J77").Value = 13098.25
$ws.Range("K77").Value = 2252.6596
$ws.Range("L77").Value = 65491.25
$ws.Range("M77").Value = 2115.3404
$ws.Range("N77").Value = -74227.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 885.83826
$ws.Range("I136").Value = 689.0172
$ws.Range("J136").Value = 2027.4
$ws.Range("K136").Value = 2067.0516
$ws.Range("L136").Value = 6082.200000000001
$ws.Range("M136").Value = 482.9484000000002
$ws.Range("N136").Value = -11182.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 987.3333
$ws.Range("I134").Value = 687.0323
$ws.Range("J134").Value = 1652.2858
$ws.Range("K134").Value = 2061.0969
$ws.Range("L134").Value = 4956.857400000001
$ws.Range("M134").Value = 473.9031
$ws.Range("N134").Value = -10026.8574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16986846
$ws.Range("I31").Value = 25001020
$ws.Range("J31").Value = 114898.42
$ws.Range("K31").Value = 25001020
$ws.Range("L31").Value = 114898.42
$ws.Range("M31").Value = -25000725
$ws.Range("N31").Value = -115488.42

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 16986846
$ws.Range("I34").Value = 25001020
$ws.Range("J34").Value = 114898.42
$ws.Range("K34").Value = 25001020
$ws.Range("L34").Value = 114898.42
$ws.Range("M34").Value = -25000818
$ws.Range("N34").Value = -115302.42

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 769.0632000000001
$ws.Range("I58").Value = 512.97015
$ws.Range("J58").Value = 1381.8572
$ws.Range("K58").Value = 512.97015
$ws.Range("L58").Value = 1381.8572
$ws.Range("M58").Value = -309.97015
$ws.Range("N58").Value = -1787.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 747.63635
$ws.Range("I94").Value = 837.3333
$ws.Range("J94").Value = 640
$ws.Range("K94").Value = 837.3333
$ws.Range("L94").Value = 640
$ws.Range("M94").Value = -386.3333
$ws.Range("N94").Value = -1542

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 11907750
$ws.Range("I132").Value = 16131849
$ws.Range("J132").Value = 3472.182
$ws.Range("K132").Value = 48395547
$ws.Range("L132").Value = 10416.546
$ws.Range("M132").Value = -48393017
$ws.Range("N132").Value = -15476.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1042.4097
$ws.Range("I134").Value = 1012.92957
$ws.Range("J134").Value = 1216.8334
$ws.Range("K134").Value = 3038.78871
$ws.Range("L134").Value = 3650.5002
$ws.Range("M134").Value = -503.7887099999998
$ws.Range("N134").Value = -8720.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 769.0632000000001
$ws.Range("I136").Value = 512.97015
$ws.Range("J136").Value = 1381.8572
$ws.Range("K136").Value = 1538.91045
$ws.Range("L136").Value = 4145.571599999999
$ws.Range("M136").Value = 1011.08955
$ws.Range("N136").Value = -9245.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 354
$ws.Range("I92").Value = 200
$ws.Range("J92").Value = 392.5
$ws.Range("K92").Value = 600
$ws.Range("L92").Value = 1177.5
$ws.Range("M92").Value = 648
$ws.Range("N92").Value = -3673.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3059.535
$ws.Range("I132").Value = 2930.7385
$ws.Range("J132").Value = 3458.1904
$ws.Range("K132").Value = 8792.2155
$ws.Range("L132").Value = 10374.5712
$ws.Range("M132").Value = -6262.2155
$ws.Range("N132").Value = -15434.5712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2938.966
$ws.Range("I136").Value = 1189.5
$ws.Range("J136").Value = 7261.1763
$ws.Range("K136").Value = 3568.5
$ws.Range("L136").Value = 21783.5289
$ws.Range("M136").Value = -1018.5
$ws.Range("N136").Value = -26883.5289

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8067352.5
$ws.Range("I132").Value = 12824102
$ws.Range("J132").Value = 1559.3043
$ws.Range("K132").Value = 38472306
$ws.Range("L132").Value = 4677.9129
$ws.Range("M132").Value = -38469776
$ws.Range("N132").Value = -9737.912899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1172.8923
$ws.Range("I136").Value = 510.07318
$ws.Range("J136").Value = 2305.2083
$ws.Range("K136").Value = 1530.21954
$ws.Range("L136").Value = 6915.624899999999
$ws.Range("M136").Value = 1019.78046
$ws.Range("N136").Value = -12015.6249
